$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.209.73'
$ws.Range('E2').Value = '  -1.55%  '
$ws.Range('D3').Value = '2.245.18'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.32'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.54%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -3.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.18'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.87%  '
$ws.Range('E11').Value = '  -2.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.17'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.103'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.851'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.92%  '
$ws.Range('D16').Value = '2.234.26'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').Value = '42.096.41'
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.40'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.47%  '
$ws.Range('E21').Value = '  +2.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +39.91%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.34%  '
$ws.Range('E27').Value = '  -1.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.67'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0821'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.25%  '
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.78'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.07%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.27'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.55%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.125'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.45'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0315'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '13.71'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.17'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.88%  '
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '62.82'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.204'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '106.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.92%  '
$ws.Range('E44').Value = '  +1.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.76%  '
$ws.Range('E46').Value = '  -0.44%  '
$ws.Range('E47').Value = '  -3.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.31'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.24'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.47%  '
$ws.Range('E51').Value = '  -0.32%  '
